# Generate Report for Archive
# - Update Status text from "Ready for handoff" to "In Translation"
#   on the Overview sheet (zh-cn / de-de columns) and on each
#   per-locale sheet's Status column.
# - Narrow the Status/locale columns to match the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.55
$overview.Columns.Item(6).ColumnWidth = 12.55

# --- zh-cn sheet ------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.55

# --- de-de sheet ------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.55
